$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2161823333333333
$ws.Range("H2").Value = 0.648547
$ws.Range("Q2").Value = 0.015495084924
$ws.Range("R2").Value = 0.139455764316
